$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-06 Tuesday" "2024-08-07 Wednesday"

Replace-Text "714÷6=119, 0" "485÷7=69, 2"
Replace-Text "137÷9=15, 2" "243÷2=121, 1"
Replace-Text "666÷2=333, 0" "558÷5=111, 3"
Replace-Text "712÷9=79, 1" "115÷4=28, 3"
Replace-Text "661÷3=220, 1" "532÷6=88, 4"

Replace-Text "675÷5=135, 0" "819÷2=409, 1"
Replace-Text "205÷4=51, 1" "361÷3=120, 1"
Replace-Text "806÷5=161, 1" "328÷7=46, 6"
Replace-Text "528÷4=132, 0" "632÷3=210, 2"
Replace-Text "258÷3=86, 0" "160÷6=26, 4"

Replace-Text "605÷5=121, 0" "679÷9=75, 4"
Replace-Text "352÷3=117, 1" "151÷2=75, 1"
Replace-Text "459÷9=51, 0" "526÷4=131, 2"
Replace-Text "880÷7=125, 5" "463÷6=77, 1"
Replace-Text "774÷4=193, 2" "687÷3=229, 0"

Replace-Text "940÷7=134, 2" "400÷3=133, 1"
Replace-Text "325÷7=46, 3" "189÷8=23, 5"
Replace-Text "962÷2=481, 0" "679÷8=84, 7"
Replace-Text "234÷9=26, 0" "294÷2=147, 0"
Replace-Text "604÷6=100, 4" "374÷4=93, 2"

Replace-Text "640÷6=106, 4" "414÷6=69, 0"
Replace-Text "845÷7=120, 5" "814÷7=116, 2"
Replace-Text "849÷4=212, 1" "825÷8=103, 1"
Replace-Text "363÷8=45, 3" "980÷6=163, 2"
Replace-Text "332÷6=55, 2" "149÷3=49, 2"
